# Add a new "RD" (hyperlink) column (H) to the players worksheet and
# populate five cells with hyperlinks to external short-links.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H
$ws.Range("H1").Value = "RD"

# Hyperlinks - added in the same order as in the original workbook so the
# generated relationship ids (rId1..rId5) line up with the target file.
$ws.Hyperlinks.Add($ws.Range("H3"), "https://t.ly/Wp6fG")
$ws.Hyperlinks.Add($ws.Range("H23"), "https://t.ly/rgqWr")
$ws.Hyperlinks.Add($ws.Range("H12"), "https://t.ly/X2Chj")

# This particular cell ends up with the plain "Hyperlink" cell style
# (no centered alignment) in the source workbook, unlike the other four.
$ws.Range("H12").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("H20"), "https://t.ly/vvdbd")
$ws.Hyperlinks.Add($ws.Range("H11"), "https://t.ly/AEuGc")

# Widen the new column to fit its contents.
$ws.Columns.Item(8).ColumnWidth = 14.5

# Match the author's final selection/view state.
[void]$ws.Range("H12").Select()
